$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("J2").Value = 4924
$ws.Range("J3").Value = 5236
$ws.Range("C4").Value = 1833
$ws.Range("F4").Value = 1894
$ws.Range("I4").Value = 1772
$ws.Range("J4").Value = 1160
$ws.Range("J5").Value = 414
$ws.Range("J6").Value = 6459
$ws.Range("C7").Value = 28377
$ws.Range("F7").Value = 24085
$ws.Range("I7").Value = 26223
$ws.Range("J7").Value = 18193

$ws = $wb.Worksheets.Item('Grant Park')
$ws.Range("J2").Value = 6
$ws.Range("J5").Value = 5
$ws.Range("J6").Value = 13

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("J3").Value = 293
$ws.Range("J6").Value = 225
$ws.Range("J7").Value = 796

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("J2").Value = 117
$ws.Range("J3").Value = 143
$ws.Range("J6").Value = 108
$ws.Range("J7").Value = 389

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("J3").Value = 92
$ws.Range("J7").Value = 273

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("J3").Value = 105
$ws.Range("J7").Value = 276

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("J3").Value = 272
$ws.Range("J6").Value = 186
$ws.Range("J7").Value = 705

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("J2").Value = 64
$ws.Range("J6").Value = 44
$ws.Range("J7").Value = 167

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("J4").Value = 5
$ws.Range("J7").Value = 150

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("J5").Value = 57
$ws.Range("J7").Value = 528
$ws.Range("J8").Value = 1160
$ws.Range("J15").Value = 195
$ws.Range("J19").Value = 526
$ws.Range("J20").Value = 384
$ws.Range("J21").Value = 47
$ws.Range("J23").Value = 176
$ws.Range("J25").Value = 90
$ws.Range("J27").Value = 100
$ws.Range("J29").Value = 1021
$ws.Range("J31").Value = 167
$ws.Range("J32").Value = 30
$ws.Range("J33").Value = 828
$ws.Range("J35").Value = 29
$ws.Range("J37").Value = 570
$ws.Range("J38").Value = 13
$ws.Range("J41").Value = 118
$ws.Range("J42").Value = 740
$ws.Range("J44").Value = 137
$ws.Range("J48").Value = 205
$ws.Range("J51").Value = 232
$ws.Range("J52").Value = 460
$ws.Range("J57").Value = 79
$ws.Range("J60").Value = 115
$ws.Range("C63").Value = 264
$ws.Range("F63").Value = 184
$ws.Range("I63").Value = 235
$ws.Range("J63").Value = 76
$ws.Range("J65").Value = 481
$ws.Range("J67").Value = 705
$ws.Range("J73").Value = 172
$ws.Range("J76").Value = 262
$ws.Range("J77").Value = 144
$ws.Range("J79").Value = 526
$ws.Range("J82").Value = 22
$ws.Range("J83").Value = 389
$ws.Range("J84").Value = 150
$ws.Range("J85").Value = 796
$ws.Range("J86").Value = 110
$ws.Range("J91").Value = 203
$ws.Range("J93").Value = 77
$ws.Range("J94").Value = 172
$ws.Range("J95").Value = 273
$ws.Range("J98").Value = 116
$ws.Range("J99").Value = 276
$ws.Range("C101").Value = 28377
$ws.Range("F101").Value = 24085
$ws.Range("I101").Value = 26223
$ws.Range("J101").Value = 18193

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("J3").Value = 199
$ws.Range("J6").Value = 163
$ws.Range("J7").Value = 570

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("J2").Value = 204
$ws.Range("J3").Value = 268
$ws.Range("J6").Value = 284
$ws.Range("J7").Value = 828

$ws = $wb.Worksheets.Item('New City')
$ws.Range("J2").Value = 136
$ws.Range("J5").Value = 13
$ws.Range("J6").Value = 172
$ws.Range("J7").Value = 481

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("J3").Value = 356
$ws.Range("J7").Value = 1021

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("J3").Value = 152
$ws.Range("J7").Value = 526

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("J6").Value = 48
$ws.Range("J7").Value = 137

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("J6").Value = 104
$ws.Range("J7").Value = 205

$ws = $wb.Worksheets.Item('River North')
$ws.Range("J3").Value = 55
$ws.Range("J6").Value = 140
$ws.Range("J7").Value = 262

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("J2").Value = 168
$ws.Range("J3").Value = 159
$ws.Range("J7").Value = 528

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range("J6").Value = 66
$ws.Range("J7").Value = 118

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("J2").Value = 161
$ws.Range("J6").Value = 380
$ws.Range("J7").Value = 740

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("J3").Value = 62
$ws.Range("J7").Value = 176

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("J3").Value = 85
$ws.Range("J7").Value = 203

$ws = $wb.Worksheets.Item('Chinatown')
$ws.Range("J6").Value = 29
$ws.Range("J7").Value = 47

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("J2").Value = 148
$ws.Range("J3").Value = 187
$ws.Range("J6").Value = 146
$ws.Range("J7").Value = 526

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("J2").Value = 104
$ws.Range("J3").Value = 129
$ws.Range("J6").Value = 103
$ws.Range("J7").Value = 384

$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range("J3").Value = 23
$ws.Range("J7").Value = 77

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("J3").Value = 33
$ws.Range("J6").Value = 94
$ws.Range("J7").Value = 172

$ws = $wb.Worksheets.Item('East Side')
$ws.Range("J2").Value = 39
$ws.Range("J7").Value = 90

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("J3").Value = 48
$ws.Range("J6").Value = 83
$ws.Range("J7").Value = 195

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("J6").Value = 68
$ws.Range("J7").Value = 116

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("J3").Value = 135
$ws.Range("J7").Value = 460

$ws = $wb.Worksheets.Item('Gold Coast')
$ws.Range("J6").Value = 18
$ws.Range("J7").Value = 29

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("J2").Value = 60
$ws.Range("J6").Value = 52
$ws.Range("J7").Value = 172

$ws = $wb.Worksheets.Item('Galewood')
$ws.Range("J2").Value = 8
$ws.Range("J7").Value = 30

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range("J3").Value = 12
$ws.Range("J7").Value = 57

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("J3").Value = 21
$ws.Range("J7").Value = 100

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("J4").Value = 58
$ws.Range("J7").Value = 110

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("J2").Value = 54
$ws.Range("J3").Value = 64
$ws.Range("J6").Value = 83
$ws.Range("J7").Value = 232

$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Range("J6").Value = 29
$ws.Range("J7").Value = 79

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("J2").Value = 40
$ws.Range("J7").Value = 115

$ws = $wb.Worksheets.Item('Sheffield & DePaul')
$ws.Range("J5").Value = 14
$ws.Range("J6").Value = 22

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("J3").Value = 52
$ws.Range("J7").Value = 144

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("J2").Value = 327
$ws.Range("J3").Value = 348
$ws.Range("J4").Value = 65
$ws.Range("J6").Value = 386
$ws.Range("J7").Value = 1160
